# #5: property aircraft done
# For the 建物 (Building) sheet, every row's property_category column (I)
# was incorrectly tagged as "land" — fix it to "building".
# For the 汽車 (Car) sheet, the single data row's property_category column
# (H) was also incorrectly tagged as "land" — fix it to "car".

$wb = $excel.ActiveWorkbook

$wsBuilding = $wb.Worksheets.Item("建物")
$wsBuilding.Range("I2:I7").Value = "building"

$wsCar = $wb.Worksheets.Item("汽車")
$wsCar.Range("H2").Value = "car"
